# Commit: "fixed test cases failing"
#
# 1. Rename the "2013 Cohort" sheet to "2013_Cohort" (space -> underscore),
#    which was presumably causing a formula/reference lookup elsewhere to
#    fail a test case that depended on an underscore-joined sheet name.
# 2. Update the sheet's remembered selection from I12 to H31 (the cell the
#    author was last working in before saving).

$wb = $excel.ActiveWorkbook

# --- 1. Rename the worksheet ---------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Name = "2013_Cohort"

# --- 2. Move the active selection to H31 ---------------------------------
$ws.Activate()
$ws.Range("H31").Select()
